# Fill in previously-missing data points in the "Performance Data" sheet.
# (Columns B = MSCI_World, C = MSCI_ACWI, D = MSCI_ACWI_IMI were sparsely
# populated; this backfills the gaps that the source feed now supplies.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 25-35: fill missing MSCI_ACWI_IMI (column D) values ---
$ws.Range("D25").Value = 629.341775
$ws.Range("D26").Value = 630.341775
$ws.Range("D27").Value = 631.341775
$ws.Range("D28").Value = 632.341775
$ws.Range("D29").Value = 633.341775
$ws.Range("D30").Value = 634.341775
$ws.Range("D31").Value = 635.341775
$ws.Range("D32").Value = 636.341775
$ws.Range("D33").Value = 637.341775
$ws.Range("D34").Value = 638.341775
$ws.Range("D35").Value = 639.341775

# --- Rows 72-87: fill missing MSCI_World (column B) values ---
$ws.Range("B72").Value = 3003.6955
$ws.Range("B73").Value = 3004.6955
$ws.Range("B74").Value = 3005.6955
$ws.Range("B75").Value = 3006.6955
$ws.Range("B76").Value = 3007.6955
$ws.Range("B77").Value = 3008.6955
$ws.Range("B78").Value = 3009.6955
$ws.Range("B79").Value = 3010.6955
$ws.Range("B80").Value = 3011.6955
$ws.Range("B81").Value = 3012.6955
$ws.Range("B82").Value = 3013.6955
$ws.Range("B83").Value = 3014.6955
$ws.Range("B84").Value = 3015.6955
$ws.Range("B85").Value = 3016.6955
$ws.Range("B86").Value = 3017.6955
$ws.Range("B87").Value = 3018.6955

# --- Rows 136-143: fill missing MSCI_World (column B) values ---
$ws.Range("B136").Value = 3662.751978
$ws.Range("B137").Value = 3663.751978
$ws.Range("B138").Value = 3664.751978
$ws.Range("B139").Value = 3665.751978
$ws.Range("B140").Value = 3666.751978
$ws.Range("B141").Value = 3667.751978
$ws.Range("B142").Value = 3668.751978
$ws.Range("B143").Value = 3669.751978

# --- Rows 138-142: also fill missing MSCI_ACWI_IMI (column D) values ---
$ws.Range("D138").Value = 832.647786
$ws.Range("D139").Value = 833.647786
$ws.Range("D140").Value = 834.647786
$ws.Range("D141").Value = 835.647786
$ws.Range("D142").Value = 836.647786

# --- Rows 188-196: fill missing MSCI_ACWI_IMI (column D) values ---
$ws.Range("D188").Value = 1292.968789
$ws.Range("D189").Value = 1293.968789
$ws.Range("D190").Value = 1294.968789
$ws.Range("D191").Value = 1295.968789
$ws.Range("D192").Value = 1296.968789
$ws.Range("D193").Value = 1297.968789
$ws.Range("D194").Value = 1298.968789
$ws.Range("D195").Value = 1299.968789
$ws.Range("D196").Value = 1300.968789

# --- Rows 268-273: fill missing MSCI_World (column B) values ---
$ws.Range("B268").Value = 11513.425105
$ws.Range("B269").Value = 11514.425105
$ws.Range("B270").Value = 11515.425105
$ws.Range("B271").Value = 11516.425105
$ws.Range("B272").Value = 11517.425105
$ws.Range("B273").Value = 11518.425105

# --- Rows 274-279: fill missing MSCI_ACWI_IMI (column D) values ---
$ws.Range("D274").Value = 2629.354585
$ws.Range("D275").Value = 2630.354585
$ws.Range("D276").Value = 2631.354585
$ws.Range("D277").Value = 2632.354585
$ws.Range("D278").Value = 2633.354585
$ws.Range("D279").Value = 2634.354585

# --- Rows 301-305: fill missing MSCI_World (column B) values ---
$ws.Range("B301").Value = 12677.970282
$ws.Range("B302").Value = 12678.970282
$ws.Range("B303").Value = 12679.970282
$ws.Range("B304").Value = 12680.970282
$ws.Range("B305").Value = 12681.970282

# --- Rows 305-321: fill missing MSCI_ACWI_IMI (column D) values ---
$ws.Range("D305").Value = 2877.617573
$ws.Range("D306").Value = 2878.617573
$ws.Range("D307").Value = 2879.617573
$ws.Range("D308").Value = 2880.617573
$ws.Range("D309").Value = 2881.617573
$ws.Range("D310").Value = 2882.617573
$ws.Range("D311").Value = 2883.617573
$ws.Range("D312").Value = 2884.617573
$ws.Range("D313").Value = 2885.617573
$ws.Range("D314").Value = 2886.617573
$ws.Range("D315").Value = 2887.617573
$ws.Range("D316").Value = 2888.617573
$ws.Range("D317").Value = 2889.617573
$ws.Range("D318").Value = 2890.617573
$ws.Range("D319").Value = 2891.617573
$ws.Range("D320").Value = 2892.617573
$ws.Range("D321").Value = 2893.617573

# --- Update the sheet view: scroll position + active selection ---
$ws.Activate()
$ws.Range("E313").Select()
$excel.ActiveWindow.ScrollRow = 286
$excel.ActiveWindow.ScrollColumn = 1
